$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# -----------------------------------------------------------------------
# The TestResult column (currently D) moves to the new column F; two new
# columns - StartDate / EndDate - take over D and E. G3 (an isolated,
# otherwise-empty styled cell) is left completely untouched.
# -----------------------------------------------------------------------

# 1) Give the new F column (TestResult) the exact look the old D column
#    header used to have, then fill in the (updated) TestResult values.
$ws.Range("D1").Copy() | Out-Null
$ws.Range("F1").PasteSpecial(-4122) | Out-Null
$ws.Range("F1").Value = "TestResult"
$ws.Range("F2").Value = "Failed"
$ws.Range("F3").ClearContents()
$ws.Range("F4").ClearContents()

# 2) Re-purpose D1/E1 as the StartDate/EndDate headers: start from the
#    plain default look (copied from A1) and add the yellow highlight.
$ws.Range("A1").Copy() | Out-Null
$ws.Range("D1").PasteSpecial(-4122) | Out-Null
$ws.Range("A1").Copy() | Out-Null
$ws.Range("E1").PasteSpecial(-4122) | Out-Null
$ws.Range("D1").Value = "StartDate"
$ws.Range("E1").Value = "EndDate"
$ws.Range("D1:E1").Interior.Color = 65535

# 3) StartDate / EndDate values (stored as plain text, matching the
#    source data - force text with a leading apostrophe so Excel does not
#    silently convert them to real date serials).
$ws.Range("D2").Value = "'01/10/2024"
$ws.Range("E2").Value = "'10/11/2024"
$ws.Range("D3").Value = "'01/10/2024"
$ws.Range("E3").Value = "'14/11/2024"
$ws.Range("D4").Value = "'01/10/2024"
$ws.Range("E4").Value = "'14/11/2024"

# Give D2 a date-look number format, then fan that exact style out to the
# rest of the new date cells so they all share one style entry.
$ws.Range("D2").NumberFormat = "mm-dd-yy"
$ws.Range("D2").Copy() | Out-Null
$ws.Range("E2").PasteSpecial(-4122) | Out-Null
$ws.Range("D3").PasteSpecial(-4122) | Out-Null
$ws.Range("E3").PasteSpecial(-4122) | Out-Null
$ws.Range("D4").PasteSpecial(-4122) | Out-Null
$ws.Range("E4").PasteSpecial(-4122) | Out-Null

# 4) Column widths for D (wider now), and the two new columns E/F.
$ws.Range("D1").ColumnWidth = 17.09
$ws.Range("E1").ColumnWidth = 14.09
$ws.Range("F1").ColumnWidth = 14

# -----------------------------------------------------------------------
# Conditional formatting: identical rules/colors, now targeting column F.
# -----------------------------------------------------------------------
$fcHeader = $ws.Range("D1").FormatConditions.Item(1)
$fcHeader.ModifyAppliesToRange($ws.Range("F3:F1048576")) | Out-Null

$fcF1a = $ws.Range("F1").FormatConditions.Add(1, 3, '"""Failed"""')
$fcF1a.Interior.Color = 192
$fcF1b = $ws.Range("F1").FormatConditions.Add(1, 3, '"Passed"')
$fcF1b.Interior.Color = 5287936

$fcD2 = $ws.Range("D2").FormatConditions.Item(1)
$fcD2.ModifyAppliesToRange($ws.Range("F2")) | Out-Null

$wb.Save()
